$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4 stays empty but must exist as a real cell (matches <c r="A4" t="inlineStr"/> in the target).
# Re-assigning Style to itself is a no-op format-wise but forces the engine to emit the cell.
$ws.Range("A4").Style = $ws.Range("A4").Style

$ws.Range("B4").Value = "What is for you Brio Maté ?"
$ws.Range("C4").Value = "Pour moi, Brio Maté est une très bonne marque de maté parce que c'est sustainable et que le goût est très bon."
$ws.Range("D4").Value = "Qu'est-ce que vous aimez le plus à propos de Brio Maté ?"
$ws.Range("E4").Value = "Concernant le Brio Maté, ce que j'aime le plus, c'est le nouveau goût avec le limon qui est vraiment très bien."
$ws.Range("F4").Value = "Quels autres avantages trouvez-vous dans le Brio Maté par rapport à d'autres marques de maté ?"
$ws.Range("G4").Value = "Les autres avantages, c'est aussi qu'ils livrent super vite. Là, j'ai tout réussi en une journée, donc c'est incroyable."
